$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 20330051920245
$ws.Range("A3").Value = 20330051920250

$ws.Range("B2").Value = "PONCE"
$ws.Range("B3").Value = "ROMAN"

$ws.Range("C2").Value = "GOMEZ"
$ws.Range("C3").Value = "ANTONIO"

$ws.Range("D2").Value = "ALETHIA LUCIA"
$ws.Range("D3").Value = "FABIOLA"

$ws.Range("E2").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS BACTERIOLÓGICAS"
$ws.Range("E3").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS BACTERIOLÓGICAS"

$ws.Range("F2").Value = "3ALCM"
$ws.Range("F3").Value = "3ALCM"

$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 6
